# Fruta / hortaliza, semanal
# Insert 3 new (more recent) weekly rows at the top of the data block that
# starts at row 176, pushing the existing rows 176:224 down to 179:227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 176; this shifts the old rows
# 176:224 down to 179:227 and grows the used range to A1:T227.
$ws.Rows("176:178").Insert()

# Common (constant) values shared by every data row in this block.
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103006
$categoria = "Nectarín"
$origen    = "Región de O'Higgins"

# New row 176: Super Queen / Especial
$ws.Range("A176").Value = $mercadoId
$ws.Range("B176").Value = $mercado
$ws.Range("C176").Value = $region
$ws.Range("D176").Value = 44559
$ws.Range("E176").Value = $codreg
$ws.Range("F176").Value = $tipo
$ws.Range("G176").Value = $prodId
$ws.Range("H176").Value = $producto
$ws.Range("I176").Value = $catId
$ws.Range("J176").Value = $categoria
$ws.Range("K176").Value = "Super Queen"
$ws.Range("L176").Value = "Especial"
$ws.Range("M176").Value = 50
$ws.Range("N176").Value = 15000
$ws.Range("O176").Value = 15000
$ws.Range("P176").Value = 15000
$ws.Range("Q176").Value = "$/caja 16 kilos empedrada"
$ws.Range("R176").Value = $origen
$ws.Range("S176").Value = 938
$ws.Range("T176").Value = 16

# New row 177: Super Queen / Primera
$ws.Range("A177").Value = $mercadoId
$ws.Range("B177").Value = $mercado
$ws.Range("C177").Value = $region
$ws.Range("D177").Value = 44559
$ws.Range("E177").Value = $codreg
$ws.Range("F177").Value = $tipo
$ws.Range("G177").Value = $prodId
$ws.Range("H177").Value = $producto
$ws.Range("I177").Value = $catId
$ws.Range("J177").Value = $categoria
$ws.Range("K177").Value = "Super Queen"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 50
$ws.Range("N177").Value = 13000
$ws.Range("O177").Value = 13000
$ws.Range("P177").Value = 13000
$ws.Range("Q177").Value = "$/caja 16 kilos empedrada"
$ws.Range("R177").Value = $origen
$ws.Range("S177").Value = 812
$ws.Range("T177").Value = 16

# New row 178: Super Queen / Segunda
$ws.Range("A178").Value = $mercadoId
$ws.Range("B178").Value = $mercado
$ws.Range("C178").Value = $region
$ws.Range("D178").Value = 44559
$ws.Range("E178").Value = $codreg
$ws.Range("F178").Value = $tipo
$ws.Range("G178").Value = $prodId
$ws.Range("H178").Value = $producto
$ws.Range("I178").Value = $catId
$ws.Range("J178").Value = $categoria
$ws.Range("K178").Value = "Super Queen"
$ws.Range("L178").Value = "Segunda"
$ws.Range("M178").Value = 50
$ws.Range("N178").Value = 11000
$ws.Range("O178").Value = 11000
$ws.Range("P178").Value = 11000
$ws.Range("Q178").Value = "$/caja 16 kilos empedrada"
$ws.Range("R178").Value = $origen
$ws.Range("S178").Value = 688
$ws.Range("T178").Value = 16
